# This edit inserts a new weekly price record for Orégano
# (Mercado Mayorista Lo Valledor de Santiago) into the existing table.
# The new record is inserted as row 271, pushing the former rows
# 271-283 down to 272-284 (the sheet grows from A1:R283 to A1:R284).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 271; this shifts all rows
# 271..283 down by one (they keep their original values/styles).
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new data point.
$ws.Range("A271").Value = 6
$ws.Range("B271").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C271").Value = "Metropolitana"
$ws.Range("D271").Value = 44939
$ws.Range("E271").Value = 13
$ws.Range("F271").Value = 100112029
$ws.Range("G271").Value = "Orégano"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 49
$ws.Range("K271").Value = 16000
$ws.Range("L271").Value = 17000
$ws.Range("M271").Value = 16469
$ws.Range("N271").Value = "$/docena de atados"
$ws.Range("O271").Value = "Región Metropolitana"
$ws.Range("P271").Value = 5490
$ws.Range("Q271").Value = 3
$ws.Range("R271").Value = "Hortaliza"
